$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row at row 6 for the Multiplan URL setting (shifts MaxRetry and
# the Aetna/Cigna/BCBS/Anthem/... Y-N provider rows down by one).
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the Multiplan URL config entry.
$ws.Range("A6").Value = "URL_Multiplan"
$ws.Range("B6").Value = "https://www.multiplan.com/webcenter/portal/ProviderSearch?_afrLoop=4494409455183375&_afrWindowMode=2&Adf-Window-Id=g6y1unwep&_afrFS=16&_afrMT=screen&_afrMFW=1366&_afrMFH=576&_afrMFDW=1366&_afrMFDH=768&_afrMFC=8&_afrMFCI=0&_afrMFM=0&_afrMFR=96&_afrMFG=0&_afrMFS=0&_afrMFO=0"

# County validation removed from Cigna -> Aetna and Cigna flags flipped to "N".
$ws.Range("B8").Value = "N"
$ws.Range("B9").Value = "N"

# Anthem selector fixed -> Anthem flag flipped to "Y".
$ws.Range("B11").Value = "Y"

# Delta Dental row now also carries the "Y or N" description, like the rest
# of the provider rows above it.
$ws.Range("C18").Value = "Y or N"

# Move the active selection to B11, matching where the author's cursor ended
# up after making the Anthem fix.
$ws.Range("B11").Select()
